# Slider battery + slider windfarm addition to the generic Batteries_ProjectTemplate
# Adds a second data row (row 2) to the "batteries" sheet describing the
# SLIDER_GB asset, and nudges a couple of cosmetic sheet-view properties
# (the A column width and the active selection) to match the edited file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# gc_id / gc_name / owner_id
$ws.Range("A2").Value = "SLIDER_GB"
$ws.Range("B2").Value = "SLIDER_GB"
$ws.Range("C2").Value = "SLIDER_GB_Owner"

# initially_active
$ws.Range("K2").Value = $false

# storage_capacity_kwh / capacity_electric_kw
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0

# connection_capacity_kw / contracted_delivery_capacity_kw / contracted_feed_in_capacity_kw
$ws.Range("N2").Value = 5000
$ws.Range("O2").Value = 5000
$ws.Range("P2").Value = 5000

# default_operation_mode
$ws.Range("Q2").Value = "NODAL_PRICING"

# latitude / longitude
$ws.Range("R2").Value = 52
$ws.Range("S2").Value = 5

# Column A widened to fit the new longer gc_id values.
$ws.Columns.Item(1).ColumnWidth = 9

# Active cell/selection moved as part of the edit.
$null = $ws.Range("C13").Select()
